$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1721.9231
$ws.Range("I70").Value = 1538.2
$ws.Range("J70").Value = 2334.3333
$ws.Range("K70").Value = 4614.6
$ws.Range("L70").Value = 7002.999899999999
$ws.Range("M70").Value = -4344.6
$ws.Range("N70").Value = -7542.999899999999

$ws.Range("H73").Value = 1721.9231
$ws.Range("I73").Value = 1538.2
$ws.Range("J73").Value = 2334.3333
$ws.Range("K73").Value = 4614.6
$ws.Range("L73").Value = 7002.999899999999
$ws.Range("M73").Value = -3678.6
$ws.Range("N73").Value = -8874.999899999999

$ws.Range("H74").Value = 7252.8887
$ws.Range("I74").Value = 7252.8887
$ws.Range("K74").Value = 7252.8887
$ws.Range("M74").Value = -6316.8887

$ws.Range("H77").Value = 7252.8887
$ws.Range("I77").Value = 7252.8887
$ws.Range("K77").Value = 36264.4435
$ws.Range("M77").Value = -31584.4435

$ws.Range("H137").Value = 12711.518
$ws.Range("I137").Value = 4970.6665
$ws.Range("J137").Value = 21005.285
$ws.Range("K137").Value = 14911.9995
$ws.Range("L137").Value = 63015.855
$ws.Range("M137").Value = -12361.9995
$ws.Range("N137").Value = -68115.855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 5398
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H45").Value = 7770.6665
$ws.Range("I45").Value = 7992.4707
$ws.Range("K45").Value = 7992.4707
$ws.Range("M45").Value = -7615.4707

$ws.Range("H61").Value = 11179.225
$ws.Range("I61").Value = 9513.643
$ws.Range("K61").Value = 9513.643
$ws.Range("M61").Value = -9301.643

$ws.Range("H74").Value = 16798.6
$ws.Range("I74").Value = 18769.75
$ws.Range("K74").Value = 18769.75
$ws.Range("M74").Value = -17895.75

$ws.Range("H77").Value = 16798.6
$ws.Range("I77").Value = 18769.75
$ws.Range("K77").Value = 93848.75
$ws.Range("M77").Value = -89480.75

$ws.Range("H132").Value = 1288.375
$ws.Range("I132").Value = 1189.2632
$ws.Range("K132").Value = 3567.7896
$ws.Range("M132").Value = -1037.7896

$ws.Range("H136").Value = 11179.225
$ws.Range("I136").Value = 9513.643
$ws.Range("K136").Value = 28540.929
$ws.Range("M136").Value = -25990.929

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3007.125
$ws.Range("I105").Value = 2727.4285
$ws.Range("K105").Value = 2727.4285
$ws.Range("M105").Value = -980.4285

$ws.Range("H134").Value = 11812.576
$ws.Range("I134").Value = 6275.2915
$ws.Range("J134").Value = 26578.666
$ws.Range("K134").Value = 18825.8745
$ws.Range("L134").Value = 79735.99800000001
$ws.Range("M134").Value = -16290.8745
$ws.Range("N134").Value = -84805.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3021.9092
$ws.Range("I31").Value = 1637.909
$ws.Range("J31").Value = 4405.909
$ws.Range("K31").Value = 1637.909
$ws.Range("L31").Value = 4405.909
$ws.Range("M31").Value = -1342.909
$ws.Range("N31").Value = -4995.909

$ws.Range("H34").Value = 3021.9092
$ws.Range("I34").Value = 1637.909
$ws.Range("J34").Value = 4405.909
$ws.Range("K34").Value = 1637.909
$ws.Range("L34").Value = 4405.909
$ws.Range("M34").Value = -1435.909
$ws.Range("N34").Value = -4809.909

$ws.Range("H62").Value = 76667.36
$ws.Range("I62").Value = 202156.6
$ws.Range("J62").Value = 6951.1113
$ws.Range("K62").Value = 202156.6
$ws.Range("L62").Value = 6951.1113
$ws.Range("M62").Value = -201532.6
$ws.Range("N62").Value = -8199.1113

$ws.Range("H65").Value = 76667.36
$ws.Range("I65").Value = 202156.6
$ws.Range("J65").Value = 6951.1113
$ws.Range("K65").Value = 1010783
$ws.Range("L65").Value = 34755.5565
$ws.Range("M65").Value = -1007663
$ws.Range("N65").Value = -40995.5565

$ws.Range("H134").Value = 3547.9016
$ws.Range("I134").Value = 2727.9285
$ws.Range("K134").Value = 8183.7855
$ws.Range("M134").Value = -5648.7855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8088.8
$ws.Range("I70").Value = 7470.0713
$ws.Range("J70").Value = 9532.5
$ws.Range("K70").Value = 7470.0713
$ws.Range("L70").Value = 9532.5
$ws.Range("M70").Value = -7200.0713
$ws.Range("N70").Value = -10072.5

$ws.Range("H73").Value = 8088.8
$ws.Range("I73").Value = 7470.0713
$ws.Range("J73").Value = 9532.5
$ws.Range("K73").Value = 7470.0713
$ws.Range("L73").Value = 9532.5
$ws.Range("M73").Value = -6534.0713
$ws.Range("N73").Value = -11404.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7734.033
$ws.Range("I16").Value = 1019.7
$ws.Range("J16").Value = 21162.7
$ws.Range("K16").Value = 1019.7
$ws.Range("L16").Value = 21162.7
$ws.Range("M16").Value = -849.7
$ws.Range("N16").Value = -21502.7

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H40").Value = 6268.273
$ws.Range("I40").Value = 6268.273
$ws.Range("K40").Value = 6268.273
$ws.Range("M40").Value = -6132.273

$ws.Range("H74").Value = 33372.223
$ws.Range("I74").Value = 24939.6
$ws.Range("J74").Value = 43913
$ws.Range("K74").Value = 24939.6
$ws.Range("L74").Value = 43913
$ws.Range("M74").Value = -23941.6
$ws.Range("N74").Value = -45909

$ws.Range("H77").Value = 33372.223
$ws.Range("I77").Value = 24939.6
$ws.Range("J77").Value = 43913
$ws.Range("K77").Value = 74818.79999999999
$ws.Range("L77").Value = 131739
$ws.Range("M77").Value = -69826.79999999999
$ws.Range("N77").Value = -141723

$ws.Range("H93").Value = 3737.0303
$ws.Range("I93").Value = 4656.9
$ws.Range("J93").Value = 2321.8462
$ws.Range("K93").Value = 4656.9
$ws.Range("L93").Value = 2321.8462
$ws.Range("M93").Value = -3408.9
$ws.Range("N93").Value = -4817.8462

$ws.Range("H99").Value = 33119.5
$ws.Range("I99").Value = 31243.889
$ws.Range("K99").Value = 31243.889
$ws.Range("M99").Value = -28248.889

$ws.Range("H122").Value = 3339.1667
$ws.Range("J122").Value = 6390.5713
$ws.Range("L122").Value = 19171.7139
$ws.Range("N122").Value = -24071.7139

$ws.Range("H132").Value = 4985.324
$ws.Range("I132").Value = 4742.7812
$ws.Range("K132").Value = 14228.3436
$ws.Range("M132").Value = -11698.3436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4702.5713
$ws.Range("I62").Value = 3695.8572
$ws.Range("J62").Value = 5709.2856
$ws.Range("K62").Value = 3695.8572
$ws.Range("L62").Value = 5709.2856
$ws.Range("M62").Value = -3071.8572
$ws.Range("N62").Value = -6957.2856

$ws.Range("H65").Value = 4702.5713
$ws.Range("I65").Value = 3695.8572
$ws.Range("J65").Value = 5709.2856
$ws.Range("K65").Value = 18479.286
$ws.Range("L65").Value = 28546.428
$ws.Range("M65").Value = -15359.286
$ws.Range("N65").Value = -34786.428

$ws.Range("H81").Value = 1574.7727
$ws.Range("I81").Value = 760.3570999999999
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 1520.7142
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -459.7141999999999
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 1574.7727
$ws.Range("I84").Value = 760.3570999999999
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 7603.571
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -2299.571
$ws.Range("N84").Value = -40608
